$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.981.42"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "3.205.88"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "601.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "153.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.204.26"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "6.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.512"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "39.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.66%  "
$ws.Range("D15").Value = "3.731.15"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "7.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.12%  "
$ws.Range("D17").Value = "66.028.85"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "3.191.33"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "511.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "15.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.740"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "8.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "15.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "84.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "9.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "3.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "2.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("E31").Value = "  +8.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "28.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "1.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "54.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "488.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.0906"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.0421"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "2.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "8.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.305"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.121"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.64%  "
$ws.Range("D44").Value = "2.954.19"
$ws.Range("E44").Value = "  -4.35%  "
$ws.Range("D45").Value = "0.0₃0647"
$ws.Range("E45").Value = "  +6.19%  "
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "28.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("E50").Value = "  +1.66%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "2.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.59%  "
